{"js": "// The document contains a single table of two-digit \u00f7 one-digit\n// division prompts (\"NN\u00f7N=\") arranged 5-per-row, with blank spacer\n// rows between each populated row. This edit replaces the prompt\n// text in-place (same cell, same run/formatting) with a new set of\n// numbers, row by row, left to right \u2014 matching the authored diff\n// exactly. One cell (0-based row 8, col 1 \u2014 \"68\u00f75=\") is intentionally\n// left set to its original value because the diff does not touch it.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Map of \"row index in the table\" -> new values for each of its 5\n// cells, left to right. Only the rows that actually contain prompts\n// are listed; the blank spacer rows are left untouched.\nconst rowUpdates = [\n  { row: 0, values: [\"53\u00f79=\", \"61\u00f75=\", \"74\u00f76=\", \"29\u00f73=\", \"18\u00f73=\"] },\n  { row: 4, values: [\"37\u00f75=\", \"19\u00f75=\", \"37\u00f79=\", \"97\u00f74=\", \"21\u00f74=\"] },\n  { row: 8, values: [\"54\u00f72=\", \"68\u00f75=\", \"18\u00f73=\", \"66\u00f73=\", \"83\u00f74=\"] },\n  { row: 12, values: [\"50\u00f72=\", \"57\u00f78=\", \"15\u00f77=\", \"59\u00f73=\", \"83\u00f75=\"] },\n  { row: 16, values: [\"10\u00f74=\", \"59\u00f73=\", \"59\u00f75=\", \"54\u00f77=\", \"94\u00f72=\"] },\n];\n\nfor (const { row, values } of rowUpdates) {\n  for (let col = 0; col < values.length; col++) {\n    const cell = table.getCell(row, col);\n    cell.value = values[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single table of two-digit \u00f7 one-digit\n# division prompts (\"NN\u00f7N=\"), 5 per populated row, with blank spacer\n# rows between each populated row. This edit replaces the prompt text\n# in-place (same cell/run/formatting) with a new set of numbers, row\n# by row, left to right \u2014 matching the authored diff exactly. Cell\n# (row 9, col 2 \u2014 \"68\u00f75=\") is intentionally left unchanged because the\n# diff does not touch it.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowUpdates = @(\n    @{ Row = 1;  Values = @(\"53\u00f79=\", \"61\u00f75=\", \"74\u00f76=\", \"29\u00f73=\", \"18\u00f73=\") },\n    @{ Row = 5;  Values = @(\"37\u00f75=\", \"19\u00f75=\", \"37\u00f79=\", \"97\u00f74=\", \"21\u00f74=\") },\n    @{ Row = 9;  Values = @(\"54\u00f72=\", \"68\u00f75=\", \"18\u00f73=\", \"66\u00f73=\", \"83\u00f74=\") },\n    @{ Row = 13; Values = @(\"50\u00f72=\", \"57\u00f78=\", \"15\u00f77=\", \"59\u00f73=\", \"83\u00f75=\") },\n    @{ Row = 17; Values = @(\"10\u00f74=\", \"59\u00f73=\", \"59\u00f75=\", \"54\u00f77=\", \"94\u00f72=\") }\n)\n\nforeach ($update in $rowUpdates) {\n    $row = $update.Row\n    $values = $update.Values\n    for ($col = 1; $col -le $values.Length; $col++) {\n        $t.Cell($row, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
